# Apply the score-sheet correction described in the commit:
# "ft && fix: removed not valid cuts && lp gap"
# - Recomputed "Puntos al finalizar la primera rueda" / "Localias faltantes" on the Equipos sheet
# - Reassigned Local/Visita pairings and Resultado scores on the Resultados sheet
$wb = $excel.ActiveWorkbook

$wsEquipos = $wb.Worksheets.Item("Equipos")
$wsEquipos.Cells.Item(2, 4).Value = 14
$wsEquipos.Cells.Item(2, 5).Value = 3
$wsEquipos.Cells.Item(3, 5).Value = 4
$wsEquipos.Cells.Item(4, 4).Value = 1
$wsEquipos.Cells.Item(4, 5).Value = 3
$wsEquipos.Cells.Item(5, 4).Value = 15
$wsEquipos.Cells.Item(5, 5).Value = 4
$wsEquipos.Cells.Item(6, 5).Value = 4
$wsEquipos.Cells.Item(7, 4).Value = 11
$wsEquipos.Cells.Item(7, 5).Value = 3
$wsEquipos.Cells.Item(8, 4).Value = 7
$wsEquipos.Cells.Item(8, 5).Value = 4
$wsEquipos.Cells.Item(9, 4).Value = 7
$wsEquipos.Cells.Item(9, 5).Value = 3

$wsResultados = $wb.Worksheets.Item("Resultados")
$wsResultados.Cells.Item(3, 3).Value = "G"
$wsResultados.Cells.Item(3, 4).Value = "H"
$wsResultados.Cells.Item(3, 5).Value = "0:2"
$wsResultados.Cells.Item(4, 3).Value = "B"
$wsResultados.Cells.Item(4, 4).Value = "C"
$wsResultados.Cells.Item(4, 5).Value = "4:5"
$wsResultados.Cells.Item(5, 3).Value = "E"
$wsResultados.Cells.Item(5, 4).Value = "A"
$wsResultados.Cells.Item(5, 5).Value = "1:1"
$wsResultados.Cells.Item(6, 3).Value = "D"
$wsResultados.Cells.Item(6, 4).Value = "F"
$wsResultados.Cells.Item(6, 5).Value = "0:1"
$wsResultados.Cells.Item(8, 3).Value = "A"
$wsResultados.Cells.Item(8, 4).Value = "G"
$wsResultados.Cells.Item(8, 5).Value = "0:0"
$wsResultados.Cells.Item(9, 3).Value = "F"
$wsResultados.Cells.Item(9, 4).Value = "B"
$wsResultados.Cells.Item(10, 3).Value = "H"
$wsResultados.Cells.Item(10, 5).Value = "2:3"
$wsResultados.Cells.Item(11, 3).Value = "C"
$wsResultados.Cells.Item(11, 5).Value = "4:2"
$wsResultados.Cells.Item(13, 3).Value = "B"
$wsResultados.Cells.Item(13, 4).Value = "D"
$wsResultados.Cells.Item(13, 5).Value = "4:2"
$wsResultados.Cells.Item(14, 5).Value = "0:0"
$wsResultados.Cells.Item(15, 3).Value = "G"
$wsResultados.Cells.Item(15, 5).Value = "6:3"
$wsResultados.Cells.Item(16, 3).Value = "A"
$wsResultados.Cells.Item(16, 4).Value = "H"
$wsResultados.Cells.Item(16, 5).Value = "0:3"
$wsResultados.Cells.Item(18, 3).Value = "H"
$wsResultados.Cells.Item(18, 4).Value = "B"
$wsResultados.Cells.Item(18, 5).Value = "2:4"
$wsResultados.Cells.Item(19, 3).Value = "C"
$wsResultados.Cells.Item(19, 4).Value = "F"
$wsResultados.Cells.Item(19, 5).Value = "2:1"
$wsResultados.Cells.Item(20, 3).Value = "E"
$wsResultados.Cells.Item(20, 4).Value = "G"
$wsResultados.Cells.Item(20, 5).Value = "4:3"
$wsResultados.Cells.Item(21, 3).Value = "D"
$wsResultados.Cells.Item(21, 4).Value = "A"
$wsResultados.Cells.Item(21, 5).Value = "1:2"
$wsResultados.Cells.Item(23, 5).Value = "2:3"
$wsResultados.Cells.Item(24, 3).Value = "G"
$wsResultados.Cells.Item(24, 4).Value = "D"
$wsResultados.Cells.Item(24, 5).Value = "6:1"
$wsResultados.Cells.Item(25, 3).Value = "H"
$wsResultados.Cells.Item(25, 4).Value = "C"
$wsResultados.Cells.Item(25, 5).Value = "2:3"
$wsResultados.Cells.Item(26, 4).Value = "E"
$wsResultados.Cells.Item(26, 5).Value = "2:2"
$wsResultados.Cells.Item(28, 3).Value = "G"
$wsResultados.Cells.Item(28, 4).Value = "F"
$wsResultados.Cells.Item(28, 5).Value = "1:3"
$wsResultados.Cells.Item(29, 3).Value = "E"
$wsResultados.Cells.Item(29, 4).Value = "H"
$wsResultados.Cells.Item(29, 5).Value = "2:2"
$wsResultados.Cells.Item(30, 3).Value = "A"
$wsResultados.Cells.Item(30, 5).Value = "2:5"
$wsResultados.Cells.Item(31, 3).Value = "D"
$wsResultados.Cells.Item(31, 4).Value = "C"
$wsResultados.Cells.Item(31, 5).Value = "1:2"
$wsResultados.Cells.Item(33, 3).Value = "C"
$wsResultados.Cells.Item(33, 5).Value = "3:4"
$wsResultados.Cells.Item(34, 3).Value = "D"
$wsResultados.Cells.Item(34, 4).Value = "E"
$wsResultados.Cells.Item(34, 5).Value = "1:1"
$wsResultados.Cells.Item(35, 3).Value = "B"
$wsResultados.Cells.Item(35, 4).Value = "G"
$wsResultados.Cells.Item(36, 3).Value = "F"
$wsResultados.Cells.Item(36, 4).Value = "H"
$wsResultados.Cells.Item(36, 5).Value = "6:2"
$wsResultados.Cells.Item(38, 3).Value = "C"
$wsResultados.Cells.Item(38, 5).Value = "2:5"
$wsResultados.Cells.Item(39, 3).Value = "A"
$wsResultados.Cells.Item(39, 4).Value = "F"
$wsResultados.Cells.Item(40, 3).Value = "E"
$wsResultados.Cells.Item(40, 4).Value = "B"
$wsResultados.Cells.Item(40, 5).Value = "2:4"
$wsResultados.Cells.Item(41, 3).Value = "H"
$wsResultados.Cells.Item(41, 4).Value = "D"
$wsResultados.Cells.Item(41, 5).Value = "0:2"
$wsResultados.Cells.Item(43, 3).Value = "G"
$wsResultados.Cells.Item(44, 4).Value = "B"
$wsResultados.Cells.Item(45, 4).Value = "C"
$wsResultados.Cells.Item(45, 5).Value = "0:0"
$wsResultados.Cells.Item(46, 3).Value = "D"
$wsResultados.Cells.Item(46, 4).Value = "E"
$wsResultados.Cells.Item(46, 5).Value = "3:0"
$wsResultados.Cells.Item(48, 3).Value = "A"
$wsResultados.Cells.Item(48, 4).Value = "E"
$wsResultados.Cells.Item(49, 3).Value = "F"
$wsResultados.Cells.Item(49, 4).Value = "D"
$wsResultados.Cells.Item(49, 5).Value = "2:1"
$wsResultados.Cells.Item(50, 3).Value = "C"
$wsResultados.Cells.Item(50, 4).Value = "H"
$wsResultados.Cells.Item(50, 5).Value = "0:2"
$wsResultados.Cells.Item(51, 3).Value = "B"
$wsResultados.Cells.Item(51, 5).Value = "1:1"
$wsResultados.Cells.Item(53, 3).Value = "H"
$wsResultados.Cells.Item(53, 4).Value = "A"
$wsResultados.Cells.Item(53, 5).Value = "1:9"
$wsResultados.Cells.Item(54, 3).Value = "D"
$wsResultados.Cells.Item(54, 4).Value = "B"
$wsResultados.Cells.Item(54, 5).Value = "3:1"
$wsResultados.Cells.Item(55, 3).Value = "G"
$wsResultados.Cells.Item(55, 4).Value = "F"
$wsResultados.Cells.Item(55, 5).Value = "1:1"
$wsResultados.Cells.Item(56, 3).Value = "E"
$wsResultados.Cells.Item(56, 4).Value = "C"
$wsResultados.Cells.Item(56, 5).Value = "2:1"
$wsResultados.Cells.Item(58, 3).Value = "A"
$wsResultados.Cells.Item(58, 4).Value = "D"
$wsResultados.Cells.Item(58, 5).Value = "5:3"
$wsResultados.Cells.Item(59, 3).Value = "F"
$wsResultados.Cells.Item(59, 4).Value = "H"
$wsResultados.Cells.Item(59, 5).Value = "1:0"
$wsResultados.Cells.Item(60, 3).Value = "B"
$wsResultados.Cells.Item(60, 4).Value = "C"
$wsResultados.Cells.Item(60, 5).Value = "4:3"
$wsResultados.Cells.Item(61, 3).Value = "G"
$wsResultados.Cells.Item(61, 4).Value = "E"
$wsResultados.Cells.Item(61, 5).Value = "1:3"
$wsResultados.Cells.Item(63, 3).Value = "C"
$wsResultados.Cells.Item(63, 4).Value = "A"
$wsResultados.Cells.Item(63, 5).Value = "1:3"
$wsResultados.Cells.Item(64, 3).Value = "E"
$wsResultados.Cells.Item(64, 4).Value = "H"
$wsResultados.Cells.Item(64, 5).Value = "0:1"
$wsResultados.Cells.Item(65, 3).Value = "B"
$wsResultados.Cells.Item(65, 4).Value = "F"
$wsResultados.Cells.Item(65, 5).Value = "2:1"
$wsResultados.Cells.Item(66, 4).Value = "G"
$wsResultados.Cells.Item(66, 5).Value = "4:0"
$wsResultados.Cells.Item(68, 3).Value = "A"
$wsResultados.Cells.Item(68, 4).Value = "B"
$wsResultados.Cells.Item(68, 5).Value = "2:3"
$wsResultados.Cells.Item(69, 5).Value = "1:3"
$wsResultados.Cells.Item(70, 3).Value = "F"
$wsResultados.Cells.Item(70, 4).Value = "E"
$wsResultados.Cells.Item(70, 5).Value = "2:1"
$wsResultados.Cells.Item(71, 4).Value = "G"
$wsResultados.Cells.Item(71, 5).Value = "1:1"
